# Updates current FFXIV market-board derived leve-profit figures across all 8 job sheets.
# Each row's H:N block (currentAveragePrice* / LevePrice* / LeveProfit*) is refreshed
# with newly computed values; a couple of rows also gain/lose a LeveProfit cell
# depending on whether NQ/HQ profit is still negative after the refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1149493.8
$ws.Range("I11").Value = 1149493.8
$ws.Range("K11").Value = 1149493.8
$ws.Range("M11").Value = -1149353.8

$ws.Range("H136").Value = 39700.332
$ws.Range("J136").Value = 39700.332
$ws.Range("L136").Value = 39700.332
$ws.Range("N136").Value = -49900.332

$ws.Range("H137").Value = 2169.06
$ws.Range("I137").Value = 1037.871
$ws.Range("J137").Value = 4014.6843
$ws.Range("K137").Value = 3113.613
$ws.Range("L137").Value = 12044.0529
$ws.Range("M137").Value = -563.6130000000003
$ws.Range("N137").Value = -17144.0529

$ws.Range("H138").Value = 2053.402
$ws.Range("I138").Value = 683.475
$ws.Range("J138").Value = 3014.7544
$ws.Range("K138").Value = 2050.425
$ws.Range("L138").Value = 9044.263199999999
$ws.Range("M138").Value = 3089.575
$ws.Range("N138").Value = -19324.2632

$ws.Range("H141").Value = 3230.0352
$ws.Range("I141").Value = 3128.3403
$ws.Range("J141").Value = 3708
$ws.Range("K141").Value = 9385.0209
$ws.Range("L141").Value = 11124
$ws.Range("M141").Value = -4205.0209
$ws.Range("N141").Value = -21484

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3728.5134
$ws.Range("I32").Value = 3464.3538
$ws.Range("J32").Value = 5636.3335
$ws.Range("K32").Value = 3464.3538
$ws.Range("L32").Value = 5636.3335
$ws.Range("M32").Value = -3177.3538
$ws.Range("N32").Value = -6210.3335

$ws.Range("H61").Value = 921.8484999999999
$ws.Range("I61").Value = 713.3333
$ws.Range("J61").Value = 1477.8889
$ws.Range("K61").Value = 713.3333
$ws.Range("L61").Value = 1477.8889
$ws.Range("M61").Value = -501.3333
$ws.Range("N61").Value = -1901.8889

$ws.Range("H74").Value = 2126.1516
$ws.Range("I74").Value = 2022.421
$ws.Range("J74").Value = 2783.111
$ws.Range("K74").Value = 2022.421
$ws.Range("L74").Value = 2783.111
$ws.Range("M74").Value = -1148.421
$ws.Range("N74").Value = -4531.111

$ws.Range("H77").Value = 2126.1516
$ws.Range("I77").Value = 2022.421
$ws.Range("J77").Value = 2783.111
$ws.Range("K77").Value = 10112.105
$ws.Range("L77").Value = 13915.555
$ws.Range("M77").Value = -5744.105
$ws.Range("N77").Value = -22651.555

$ws.Range("H132").Value = 3289.9
$ws.Range("I132").Value = 2070.5881
$ws.Range("J132").Value = 4884.385
$ws.Range("K132").Value = 6211.7643
$ws.Range("L132").Value = 14653.155
$ws.Range("M132").Value = -3681.7643
$ws.Range("N132").Value = -19713.155

$ws.Range("H136").Value = 921.8484999999999
$ws.Range("I136").Value = 713.3333
$ws.Range("J136").Value = 1477.8889
$ws.Range("K136").Value = 2139.9999
$ws.Range("L136").Value = 4433.6667
$ws.Range("M136").Value = 410.0001000000002
$ws.Range("N136").Value = -9533.6667

$ws.Range("H137").Value = 39774
$ws.Range("J137").Value = 39774
$ws.Range("L137").Value = 39774
$ws.Range("N137").Value = -49974

$ws.Range("H138").Value = 78369.5
$ws.Range("J138").Value = 78369.5
$ws.Range("L138").Value = 78369.5
$ws.Range("N138").Value = -88649.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 31166.666
$ws.Range("J59").Value = 31166.666
$ws.Range("L59").Value = 31166.666
$ws.Range("N59").Value = -32860.666

$ws.Range("H134").Value = 1779.0758
$ws.Range("I134").Value = 982.2982
$ws.Range("J134").Value = 6825.3335
$ws.Range("K134").Value = 2946.8946
$ws.Range("L134").Value = 20476.0005
$ws.Range("M134").Value = -411.8945999999996
$ws.Range("N134").Value = -25546.0005

$ws.Range("H137").Value = 37236.668
$ws.Range("J137").Value = 37236.668
$ws.Range("L137").Value = 37236.668
$ws.Range("N137").Value = -47436.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6581277.5
$ws.Range("I31").Value = 1276.2979
$ws.Range("J31").Value = 17245418
$ws.Range("K31").Value = 1276.2979
$ws.Range("L31").Value = 17245418
$ws.Range("M31").Value = -981.2979
$ws.Range("N31").Value = -17246008

$ws.Range("H34").Value = 6581277.5
$ws.Range("I34").Value = 1276.2979
$ws.Range("J34").Value = 17245418
$ws.Range("K34").Value = 1276.2979
$ws.Range("L34").Value = 17245418
$ws.Range("M34").Value = -1074.2979
$ws.Range("N34").Value = -17245822

$ws.Range("H58").Value = 1182.48
$ws.Range("I58").Value = 1255.7902
$ws.Range("J58").Value = 869.9474
$ws.Range("K58").Value = 1255.7902
$ws.Range("L58").Value = 869.9474
$ws.Range("M58").Value = -1052.7902
$ws.Range("N58").Value = -1275.9474

$ws.Range("H104").Value = 27642.5
$ws.Range("J104").Value = 34285
$ws.Range("L104").Value = 34285
$ws.Range("N104").Value = -39527

$ws.Range("H132").Value = 2443.3142
$ws.Range("I132").Value = 2040.0667
$ws.Range("J132").Value = 4862.8
$ws.Range("K132").Value = 6120.2001
$ws.Range("L132").Value = 14588.4
$ws.Range("M132").Value = -3590.2001
$ws.Range("N132").Value = -19648.4

$ws.Range("H134").Value = 2985.2742
$ws.Range("I134").Value = 3849.2727
$ws.Range("J134").Value = 2002.1034
$ws.Range("K134").Value = 11547.8181
$ws.Range("L134").Value = 6006.3102
$ws.Range("M134").Value = -9012.8181
$ws.Range("N134").Value = -11076.3102

$ws.Range("H136").Value = 1182.48
$ws.Range("I136").Value = 1255.7902
$ws.Range("J136").Value = 869.9474
$ws.Range("K136").Value = 3767.3706
$ws.Range("L136").Value = 2609.8422
$ws.Range("M136").Value = -1217.3706
$ws.Range("N136").Value = -7709.8422

$ws.Range("H138").Value = 47143.332
$ws.Range("J138").Value = 47143.332
$ws.Range("L138").Value = 47143.332
$ws.Range("N138").Value = -57423.332

$ws.Range("H140").Value = 75683.125
$ws.Range("J140").Value = 75683.125
$ws.Range("L140").Value = 75683.125
$ws.Range("N140").Value = -86043.125

$ws.Range("H141").Value = 34850
$ws.Range("J141").Value = 34850
$ws.Range("L141").Value = 34850
$ws.Range("N141").Value = -45210

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 165.10527
$ws.Range("I23").Value = 78.625
$ws.Range("J23").Value = 228
$ws.Range("K23").Value = 235.875
$ws.Range("L23").Value = 684
$ws.Range("M23").Value = -0.875
$ws.Range("N23").Value = -1154

$ws.Range("H114").Value = 3845.2666
$ws.Range("I114").Value = 99.666664
$ws.Range("J114").Value = 4781.6665
$ws.Range("K114").Value = 298.999992
$ws.Range("L114").Value = 14344.9995
$ws.Range("M114").Value = 2955.000008
$ws.Range("N114").Value = -20852.9995

$ws.Range("H131").Value = 776.53534
$ws.Range("I131").Value = 326.4
$ws.Range("J131").Value = 827.1123700000001
$ws.Range("K131").Value = 979.1999999999999
$ws.Range("L131").Value = 2481.33711
$ws.Range("M131").Value = 4060.8
$ws.Range("N131").Value = -12561.33711

$ws.Range("H132").Value = 1473.1305
$ws.Range("I132").Value = 677.3570999999999
$ws.Range("K132").Value = 6096.2139
$ws.Range("M132").Value = -3566.2139

$ws.Range("H133").Value = 3260.25
$ws.Range("I133").Value = 4178.125
$ws.Range("K133").Value = 12534.375
$ws.Range("M133").Value = -7474.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 27420.75
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 30909.428
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 30909.428
$ws.Range("M46").Value = -2844
$ws.Range("N46").Value = -31221.428

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H80").Value = 41669120
$ws.Range("I80").Value = 83335304
$ws.Range("J80").Value = 2933.3333
$ws.Range("K80").Value = 83335304
$ws.Range("L80").Value = 2933.3333
$ws.Range("M80").Value = -83334306
$ws.Range("N80").Value = -4929.3333

$ws.Range("H83").Value = 41669120
$ws.Range("I83").Value = 83335304
$ws.Range("J83").Value = 2933.3333
$ws.Range("K83").Value = 416676520
$ws.Range("L83").Value = 14666.6665
$ws.Range("M83").Value = -416671528
$ws.Range("N83").Value = -24650.6665

$ws.Range("H132").Value = 3676.1177
$ws.Range("I132").Value = 2209.2727
$ws.Range("J132").Value = 6365.3335
$ws.Range("K132").Value = 6627.8181
$ws.Range("L132").Value = 19096.0005
$ws.Range("M132").Value = -4097.8181
$ws.Range("N132").Value = -24156.0005

$ws.Range("H137").Value = 42726.668
$ws.Range("J137").Value = 42726.668
$ws.Range("L137").Value = 42726.668
$ws.Range("N137").Value = -52926.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2589.2
$ws.Range("I46").Value = 2772.5
$ws.Range("K46").Value = 2772.5
$ws.Range("M46").Value = -2584.5

$ws.Range("H122").Value = 8433.888999999999
$ws.Range("I122").Value = 4900
$ws.Range("K122").Value = 14700
$ws.Range("M122").Value = -12250

$ws.Range("H132").Value = 47114.285
$ws.Range("I132").Value = 240000
$ws.Range("J132").Value = 14966.667
$ws.Range("K132").Value = 720000
$ws.Range("L132").Value = 44900.001
$ws.Range("M132").Value = -717470
$ws.Range("N132").Value = -49960.001

$ws.Range("H136").Value = 1987.4576
$ws.Range("I136").Value = 1117.9166
$ws.Range("J136").Value = 5781.8184
$ws.Range("K136").Value = 3353.7498
$ws.Range("L136").Value = 17345.4552
$ws.Range("M136").Value = -803.7498000000001
$ws.Range("N136").Value = -22445.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 55561056
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 55561056
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 166683168
$ws.Range("N132").Value = -166688228
$ws.Range("M132").ClearContents()
